$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking "Price" column cells (column D) must be written as TEXT,
# matching the source workbook which stores them as inline strings.
# Setting NumberFormat to "@" (Text) before the write stops Excel from
# coercing the string into a number; resetting Style to "Normal" afterwards
# removes the temporary number-format so the cell keeps its original (default) style.
$priceUpdates = @(
    @('D2', '243.64'),
    @('D3', '23.53'),
    @('D4', '5.300'),
    @('D5', '0.05769'),
    @('D6', '6.468'),
    @('D7', '3.332'),
    @('D8', '0.8144'),
    @('D9', '0.8789'),
    @('D10', '0.1380'),
    @('D11', '0.07319'),
    @('D12', '0.03093'),
    @('D13', '0.03059'),
    @('D14', '0.09327'),
    @('D15', '3.851'),
    @('D16', '0.001539'),
    @('D17', '0.04711'),
    @('D18', '0.0006018'),
    @('D19', '0.006008'),
    @('D20', '0.001294'),
    @('D21', '0.004607'),
    @('D22', '0.00008806'),
    @('D24', '2.144'),
    @('D40', '0.03767'),
    @('D41', '0.006355'),
    @('D42', '0.1054'),
    @('D43', '0.002602'),
    @('D44', '0.007192'),
    @('D45', '0.00005483'),
    @('D47', '0.6008'),
    @('D48', '0.001862')
)

foreach ($update in $priceUpdates) {
    $cell = $ws.Range($update[0])
    $cell.NumberFormat = "@"
    $cell.Value = $update[1]
    $cell.Style = "Normal"
}

# Plain text cells (coin name / link / rank-label) can be assigned directly.
$textUpdates = @(
    @('B41', 'KickToken'),
    @('C41', 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'),
    @('E41', '40KickTokenKICK'),
    @('B42', 'BKEXToken'),
    @('C42', 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'),
    @('E42', '41BKEXTokenBKK'),
    @('B43', 'CEJI'),
    @('C43', 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'),
    @('E43', '42CEJICEJIWorstin24h')
)

foreach ($update in $textUpdates) {
    $ws.Range($update[0]).Value = $update[1]
}
